# Updated cryptos list on Sun May 19 07:11:12 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21 - now Uniswap (was Polygon)
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"

# Row 22 - now Polygon (was Uniswap)
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"

# --- Price (column D) updates ---
$ws.Range("D2").Value = "67.067.68"
$ws.Range("D3").Value = "3.110.54"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.95"
$ws.Range("D5").ClearFormats()
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.33"
$ws.Range("D6").ClearFormats()
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "6.51"
$ws.Range("D9").ClearFormats()
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("D10").ClearFormats()
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.478"
$ws.Range("D11").ClearFormats()
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000248"
$ws.Range("D12").ClearFormats()
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "36.83"
$ws.Range("D13").ClearFormats()
$ws.Range("D15").Value = "3.623.19"
$ws.Range("D16").Value = "67.036.35"
$ws.Range("D18").Value = "3.109.67"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.45"
$ws.Range("D19").ClearFormats()
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "492.10"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.704"
$ws.Range("D22").ClearFormats()
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "83.86"
$ws.Range("D23").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.13"
$ws.Range("D24").ClearFormats()
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.29"
$ws.Range("D25").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.45"
$ws.Range("D26").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "28.39"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").Value = "0.0₃0947"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.999"
$ws.Range("D34").ClearFormats()
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.85"
$ws.Range("D35").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.972"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "47.09"
$ws.Range("D37").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.46"
$ws.Range("D41").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "384.65"
$ws.Range("D42").ClearFormats()
$ws.Range("D43").Value = "2.803.11"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.59"
$ws.Range("D44").ClearFormats()
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0351"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "135.59"
$ws.Range("D46").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.85"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.19"
$ws.Range("D49").ClearFormats()

# --- Volume(1h) (column E) updates ---
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("E3").Value = "  +0.40%  "
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  -0.46%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("E11").Value = "  -0.56%  "
$ws.Range("E12").Value = "  -0.33%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("E14").Value = "  -1.61%  "
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("E17").Value = "  -1.54%  "
$ws.Range("E18").Value = "  +0.32%  "
$ws.Range("E19").Value = "  +0.88%  "
$ws.Range("E20").Value = "  +1.91%  "
$ws.Range("E21").Value = "  +4.74%  "
$ws.Range("E22").Value = "  -1.74%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  -0.89%  "
$ws.Range("E25").Value = "  -3.19%  "
$ws.Range("E26").Value = "  +4.32%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("E35").Value = "  -0.59%  "
$ws.Range("E36").Value = "  -1.49%  "
$ws.Range("E37").Value = "  -2.02%  "
$ws.Range("E38").Value = "  -4.02%  "
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  +0.72%  "
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("E43").Value = "  -1.50%  "
$ws.Range("E44").Value = "  -7.83%  "
$ws.Range("E45").Value = "  -2.59%  "
$ws.Range("E46").Value = "  -0.13%  "
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("E49").Value = "  -1.71%  "
$ws.Range("E50").Value = "  -1.13%  "
$ws.Range("E51").Value = "  -2.02%  "
